$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 110076352
$ws.Range("K2").Value = ""
$ws.Range("P2").Value = "Slåttrösthammaren, Hjd"
$ws.Range("Q2").Value = 428398.2475711301
$ws.Range("R2").Value = 6968218.251045209
$ws.Range("Y2").Value = "2023-06-13"
$ws.Range("AA2").Value = "2023-06-13"
$ws.Range("AW2").Value = "Erland Lindblad"
$ws.Range("AX2").Value = "Erland Lindblad"
$ws.Range("AY2").Value = ""
$ws.Range("A3").Value = 110076358
$ws.Range("K3").Value = ""
$ws.Range("P3").Value = "Slåttrösthammaren, Hjd"
$ws.Range("Q3").Value = 428398.2475711301
$ws.Range("R3").Value = 6968218.251045209
$ws.Range("Y3").Value = "2023-06-13"
$ws.Range("AA3").Value = "2023-06-13"
$ws.Range("AW3").Value = "Erland Lindblad"
$ws.Range("AX3").Value = "Erland Lindblad"
$ws.Range("AY3").Value = ""
$ws.Range("A4").Value = 6794618
$ws.Range("B4").Value = 78072
$ws.Range("E4").Value = 229821
$ws.Range("F4").Value = "Vedflamlav"
$ws.Range("G4").Value = "Ramboldia elabens"
$ws.Range("H4").Value = "(Fr.) Kantvilas & Elix"
$ws.Range("K4").ClearContents()
$ws.Range("P4").Value = "Öster om Stor-Vävelsjön, Hjd"
$ws.Range("Q4").Value = 428683.8069340216
$ws.Range("R4").Value = 6968384.788107738
$ws.Range("Y4").Value = "2013-05-31"
$ws.Range("AA4").Value = "2013-05-31"
$ws.Range("AW4").Value = "Hugo Ström"
$ws.Range("AX4").Value = "Hugo Ström"
$ws.Range("AY4").Value = "SCA Skog Naturvärdesinventering"
$ws.Range("A5").Value = 6794619
$ws.Range("B5").Value = 78098
$ws.Range("E5").Value = 6453
$ws.Range("F5").Value = "Vedskivlav"
$ws.Range("G5").Value = "Hertelidea botryosa"
$ws.Range("H5").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("K5").ClearContents()
$ws.Range("P5").Value = "Öster om Stor-Vävelsjön, Hjd"
$ws.Range("Q5").Value = 428687.9835908828
$ws.Range("R5").Value = 6968387.439416924
$ws.Range("Y5").Value = "2013-05-31"
$ws.Range("AA5").Value = "2013-05-31"
$ws.Range("AW5").Value = "Hugo Ström"
$ws.Range("AX5").Value = "Hugo Ström"
$ws.Range("AY5").Value = "SCA Skog Naturvärdesinventering"
$ws.Range("A6").Value = 110075547
$ws.Range("B6").Value = 78098
$ws.Range("E6").Value = 6453
$ws.Range("F6").Value = "Vedskivlav"
$ws.Range("G6").Value = "Hertelidea botryosa"
$ws.Range("H6").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q6").Value = 428529.5896594818
$ws.Range("R6").Value = 6968112.548214767
$ws.Range("A7").Value = 110075549
$ws.Range("B7").Value = 78072
$ws.Range("E7").Value = 229821
$ws.Range("F7").Value = "Vedflamlav"
$ws.Range("G7").Value = "Ramboldia elabens"
$ws.Range("H7").Value = "(Fr.) Kantvilas & Elix"
$ws.Range("Q7").Value = 428529.5896594818
$ws.Range("R7").Value = 6968112.548214767
$ws.Range("A8").Value = 110076486
$ws.Range("B8").Value = 78072
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 229821
$ws.Range("F8").Value = "Vedflamlav"
$ws.Range("G8").Value = "Ramboldia elabens"
$ws.Range("H8").Value = "(Fr.) Kantvilas & Elix"
$ws.Range("Q8").Value = 428437.5493903486
$ws.Range("R8").Value = 6968005.772483499
$ws.Range("A9").Value = 110076483
$ws.Range("Q9").Value = 428437.5493903486
$ws.Range("R9").Value = 6968005.772483499
$ws.Range("A10").Value = 110075736
$ws.Range("B10").Value = 95525
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 221941
$ws.Range("F10").Value = "Plattlummer"
$ws.Range("G10").Value = "Lycopodium complanatum"
$ws.Range("H10").Value = "L."
$ws.Range("Q10").Value = 428653.7851904702
$ws.Range("R10").Value = 6968309.117210778
$ws.Range("A11").Value = 110075810
$ws.Range("B11").Value = 77177
$ws.Range("E11").Value = 353
$ws.Range("F11").Value = "Dvärgbägarlav"
$ws.Range("G11").Value = "Cladonia parasitica"
$ws.Range("H11").Value = "(Hoffm.) Hoffm."
$ws.Range("Q11").Value = 428681.7797006101
$ws.Range("R11").Value = 6968312.621547545
$ws.Range("A12").Value = 110076489
$ws.Range("Q12").Value = 428437.5493903486
$ws.Range("R12").Value = 6968005.772483499
